$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.312.06"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.681.17"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'218.77"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "'0.5275"
$ws.Range("E6").Value = "  +3.19%  "
$ws.Range("D7").Value = "'1.007"
$ws.Range("D8").Value = "'0.2709"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("D9").Value = "'0.06441"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("D10").Value = "'22.06"
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("D11").Value = "'0.07514"
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").Value = "1.696.09"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").Value = "'4.553"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "'0.5814"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "'0.000008507"
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").Value = "'64.40"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").Value = "26.341.01"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "'4.939"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "'190.09"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "'6.213"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'1.008"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'145.16"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").Value = "'7.771"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").Value = "'0.1246"
$ws.Range("E26").Value = "  +5.80%  "
$ws.Range("D27").Value = "'15.83"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").Value = "'0.06642"
$ws.Range("E28").Value = "  +10.96%  "
$ws.Range("D29").Value = "'1.359"
$ws.Range("E29").Value = "  +5.65%  "
$ws.Range("D30").Value = "'1.329"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "'3.590"
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("D32").Value = "'3.578"
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("D34").Value = "'1.028"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").Value = "'0.6217"
$ws.Range("E35").Value = "  +3.09%  "
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("E37").Value = "  +2.10%  "
$ws.Range("D38").Value = "'6.415"
$ws.Range("E38").Value = "  +5.25%  "
$ws.Range("D39").Value = "'0.01623"
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").Value = "1.108.42"
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("D41").Value = "'0.8780"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").Value = "'1.015"
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").Value = "'100.74"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").Value = "1.846.68"
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("D45").Value = "'0.00000000110"
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("D46").Value = "'56.94"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("D47").Value = "'1.011"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "'8.133"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "'0.4300"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "'6.072"
$ws.Range("E51").Value = "  +3.05%  "
